# This workbook tracks Bahamut-server crafting-Leve profitability per job.
# Columns H:N (currentAveragePrice[NQ/HQ], LevePrice[NQ/HQ], LeveProfit[NQ/HQ])
# are refreshed from the latest Universalis market snapshot for the Leve rows
# below, one table per job sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H80").Value = 599110.5600000001
$ws.Range("I80").Value = 882.5
$ws.Range("J80").Value = 1556275.4
$ws.Range("K80").Value = 2647.5
$ws.Range("L80").Value = 4668826.199999999
$ws.Range("M80").Value = -1649.5
$ws.Range("N80").Value = -4670822.199999999
$ws.Range("H83").Value = 599110.5600000001
$ws.Range("I83").Value = 882.5
$ws.Range("J83").Value = 1556275.4
$ws.Range("K83").Value = 7942.5
$ws.Range("L83").Value = 14006478.6
$ws.Range("M83").Value = -2950.5
$ws.Range("N83").Value = -14016462.6
$ws.Range("H92").Value = 2783.4736
$ws.Range("I92").Value = 3178.7334
$ws.Range("K92").Value = 3178.7334
$ws.Range("M92").Value = -1930.7334
$ws.Range("H98").Value = 2237.1724
$ws.Range("I98").Value = 2175.8096
$ws.Range("J98").Value = 2398.25
$ws.Range("K98").Value = 2175.8096
$ws.Range("L98").Value = 2398.25
$ws.Range("M98").Value = -677.8096
$ws.Range("N98").Value = -5394.25
$ws.Range("H122").Value = 2237.1724
$ws.Range("I122").Value = 2175.8096
$ws.Range("J122").Value = 2398.25
$ws.Range("K122").Value = 6527.4288
$ws.Range("L122").Value = 7194.75
$ws.Range("M122").Value = -4077.4288
$ws.Range("N122").Value = -12094.75
$ws.Range("H138").Value = 4365.959
$ws.Range("I138").Value = 2654.8235
$ws.Range("J138").Value = 5275
$ws.Range("K138").Value = 7964.470499999999
$ws.Range("L138").Value = 15825
$ws.Range("M138").Value = -2824.470499999999
$ws.Range("N138").Value = -26105

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H32").Value = 20451.688
$ws.Range("I32").Value = 14279.472
$ws.Range("J32").Value = 53686.69
$ws.Range("K32").Value = 14279.472
$ws.Range("L32").Value = 53686.69
$ws.Range("M32").Value = -13992.472
$ws.Range("N32").Value = -54260.69
$ws.Range("H74").Value = 1103.9788
$ws.Range("I74").Value = 1050.6923
$ws.Range("J74").Value = 1363.75
$ws.Range("K74").Value = 1050.6923
$ws.Range("L74").Value = 1363.75
$ws.Range("M74").Value = -176.6922999999999
$ws.Range("N74").Value = -3111.75
$ws.Range("H77").Value = 1103.9788
$ws.Range("I77").Value = 1050.6923
$ws.Range("J77").Value = 1363.75
$ws.Range("K77").Value = 5253.461499999999
$ws.Range("L77").Value = 6818.75
$ws.Range("M77").Value = -885.4614999999994
$ws.Range("N77").Value = -15554.75
$ws.Range("H122").Value = 1571.7778
$ws.Range("I122").Value = 1576
$ws.Range("K122").Value = 4728
$ws.Range("M122").Value = -2278
$ws.Range("H132").Value = 2303.7083
$ws.Range("I132").Value = 1868.1052
$ws.Range("J132").Value = 3959
$ws.Range("K132").Value = 5604.3156
$ws.Range("L132").Value = 11877
$ws.Range("M132").Value = -3074.3156
$ws.Range("N132").Value = -16937

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H86").Value = 8125.2915
$ws.Range("I86").Value = 12718.091
$ws.Range("K86").Value = 12718.091
$ws.Range("M86").Value = -11595.091
$ws.Range("H89").Value = 8125.2915
$ws.Range("I89").Value = 12718.091
$ws.Range("K89").Value = 63590.455
$ws.Range("M89").Value = -57974.455
$ws.Range("H94").Value = 630.6799999999999
$ws.Range("I94").Value = 535.6316
$ws.Range("J94").Value = 931.6667
$ws.Range("K94").Value = 535.6316
$ws.Range("L94").Value = 931.6667
$ws.Range("M94").Value = -84.63160000000005
$ws.Range("N94").Value = -1833.6667
$ws.Range("H134").Value = 80932.03999999999
$ws.Range("I134").Value = 4116.7896
$ws.Range("J134").Value = 289430.56
$ws.Range("K134").Value = 12350.3688
$ws.Range("L134").Value = 868291.6799999999
$ws.Range("M134").Value = -9815.3688
$ws.Range("N134").Value = -873361.6799999999

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 3180.6924
$ws.Range("I31").Value = 2996.1482
$ws.Range("J31").Value = 3595.9167
$ws.Range("K31").Value = 2996.1482
$ws.Range("L31").Value = 3595.9167
$ws.Range("M31").Value = -2701.1482
$ws.Range("N31").Value = -4185.9167
$ws.Range("H34").Value = 3180.6924
$ws.Range("I34").Value = 2996.1482
$ws.Range("J34").Value = 3595.9167
$ws.Range("K34").Value = 2996.1482
$ws.Range("L34").Value = 3595.9167
$ws.Range("M34").Value = -2794.1482
$ws.Range("N34").Value = -3999.9167
$ws.Range("H86").Value = 4069.5908
$ws.Range("I86").Value = 1443.3077
$ws.Range("K86").Value = 1443.3077
$ws.Range("M86").Value = -320.3077000000001
$ws.Range("H89").Value = 4069.5908
$ws.Range("I89").Value = 1443.3077
$ws.Range("K89").Value = 7216.538500000001
$ws.Range("M89").Value = -1600.538500000001
$ws.Range("H122").Value = 1635.1666
$ws.Range("I122").Value = 1866.5
$ws.Range("J122").Value = 1172.5
$ws.Range("K122").Value = 5599.5
$ws.Range("L122").Value = 3517.5
$ws.Range("M122").Value = -3149.5
$ws.Range("N122").Value = -8417.5
$ws.Range("H132").Value = 2018.6666
$ws.Range("I132").Value = 1422.5
$ws.Range("K132").Value = 4267.5
$ws.Range("M132").Value = -1737.5

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H80").Value = 8549.875
$ws.Range("J80").Value = 8666.666999999999
$ws.Range("L80").Value = 26000.001
$ws.Range("N80").Value = -27872.001
$ws.Range("H83").Value = 8549.875
$ws.Range("J83").Value = 8666.666999999999
$ws.Range("L83").Value = 78000.003
$ws.Range("N83").Value = -87360.003
$ws.Range("H131").Value = 41835896
$ws.Range("I131").Value = 100202260
$ws.Range("J131").Value = 145633
$ws.Range("K131").Value = 300606780
$ws.Range("L131").Value = 436899
$ws.Range("M131").Value = -300601740
$ws.Range("N131").Value = -446979

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H102").Value = 1942.4
$ws.Range("I102").Value = 1178
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1178
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = 444
$ws.Range("N102").Value = -8244
$ws.Range("H132").Value = 3051.8223
$ws.Range("I132").Value = 2694.861
$ws.Range("J132").Value = 4479.6665
$ws.Range("K132").Value = 8084.583
$ws.Range("L132").Value = 13438.9995
$ws.Range("M132").Value = -5554.583
$ws.Range("N132").Value = -18498.9995

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H46").Value = 1580
$ws.Range("I46").Value = 1580
$ws.Range("K46").Value = 1580
$ws.Range("M46").Value = -1392
$ws.Range("H82").Value = 2795.05
$ws.Range("I82").Value = 2849.9167
$ws.Range("J82").Value = 2712.75
$ws.Range("K82").Value = 2849.9167
$ws.Range("L82").Value = 2712.75
$ws.Range("M82").Value = -2488.9167
$ws.Range("N82").Value = -3434.75
$ws.Range("H85").Value = 2795.05
$ws.Range("I85").Value = 2849.9167
$ws.Range("J85").Value = 2712.75
$ws.Range("K85").Value = 2849.9167
$ws.Range("L85").Value = 2712.75
$ws.Range("M85").Value = -1601.9167
$ws.Range("N85").Value = -5208.75

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H122").Value = 1677.1428
$ws.Range("I122").Value = 1448
$ws.Range("K122").Value = 4344
$ws.Range("M122").Value = -1894
$ws.Range("H126").Value = 1000.7778
$ws.Range("I126").Value = 1000.7778
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3002.3334
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -532.3334
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 4915.357
$ws.Range("I132").Value = 1396.3
$ws.Range("J132").Value = 13713
$ws.Range("K132").Value = 4188.9
$ws.Range("L132").Value = 41139
$ws.Range("M132").Value = -1658.9
$ws.Range("N132").Value = -46199
$ws.Range("H136").Value = 3585.0286
$ws.Range("I136").Value = 645.92
$ws.Range("J136").Value = 10932.8
$ws.Range("K136").Value = 1937.76
$ws.Range("L136").Value = 32798.39999999999
$ws.Range("M136").Value = 612.2400000000002
$ws.Range("N136").Value = -37898.39999999999
